{"js": "// Cover letter for initial submission\n// 1) Date: \"November 13\" -> \"December 2\"\n// 2) Salutation: \"Casadevall\" -> \"Imperiale\"\n// 3) Journal: \"mBio\" -> \"mSphere\" (italic), split across two runs as in the\n//    authored edit (\"m\" then \"Sphere\")\n\nconst body = context.document.body;\n\n// 1) Update the date.\nconst dateResults = body.search(\"November 13\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"December 2\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Update the addressee's surname.\nconst nameResults = body.search(\"Casadevall\", { matchCase: true, matchWholeWord: false });\nnameResults.load(\"items\");\nawait context.sync();\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Imperiale\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Update the target journal name (kept italic).\nconst journalResults = body.search(\"mBio\", { matchCase: true, matchWholeWord: false });\njournalResults.load(\"items\");\nawait context.sync();\nif (journalResults.items.length > 0) {\n  const journalRange = journalResults.items[0];\n  journalRange.insertText(\"m\", Word.InsertLocation.replace);\n  const sphereRange = journalRange.insertText(\"Sphere\", Word.InsertLocation.after);\n  sphereRange.font.set({ italic: true });\n}\nawait context.sync();\n", "ps1": "# Cover letter for initial submission\n# 1) Date: \"November 13\" -> \"December 2\"\n# 2) Salutation: \"Casadevall\" -> \"Imperiale\"\n# 3) Journal: \"mBio\" -> \"mSphere\" (formatting, e.g. italics, is preserved\n#    automatically since Find/Replace edits the matched run in place)\n\n$d = $word.ActiveDocument\n\n# 1) Update the date.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"November 13\", $false, $false, $false, $false, $false, $true, 1, $false, \"December 2\", 2)\n\n# 2) Update the addressee's surname.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Casadevall\", $false, $false, $false, $false, $false, $true, 1, $false, \"Imperiale\", 2)\n\n# 3) Update the target journal name.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"mBio\", $false, $false, $false, $false, $false, $true, 1, $false, \"mSphere\", 2)\n"}
